$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.951.32"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -5.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.882.83"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.58%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "545.59"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "124.43"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.93%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.89%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.880.86"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.45%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -8.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.67"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -7.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.431"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000210"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -5.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.86"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.00%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.357.97"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.881.55"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.37%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +5.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "56.928.29"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -5.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "401.73"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -6.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.72"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.662"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.76"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -5.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.58"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "77.43"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.53%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.44"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.16"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.55%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "24.46"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.90"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0972"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.906"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.39"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.00"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -11.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "47.79"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.13"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0627"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.105"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0334"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -6.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.42"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.65%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "360.09"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.26%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.596.54"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.44%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "119.80"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.69%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.85%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.86%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.90%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.39"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.93%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.17%  "
